$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1771.1515
$ws.Range("I40").Value = 1505.1578
$ws.Range("J40").Value = 2132.1428
$ws.Range("K40").Value = 1505.1578
$ws.Range("L40").Value = 2132.1428
$ws.Range("M40").Value = -1330.1578
$ws.Range("N40").Value = -2482.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3758.3333
$ws.Range("I74").Value = 3700
$ws.Range("J74").Value = 4050
$ws.Range("K74").Value = 3700
$ws.Range("L74").Value = 4050
$ws.Range("M74").Value = -2764
$ws.Range("N74").Value = -5922

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3048.4849
$ws.Range("I76").Value = 3027.5862
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3027.5862
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -2712.5862
$ws.Range("N76").Value = -3830

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3758.3333
$ws.Range("I77").Value = 3700
$ws.Range("J77").Value = 4050
$ws.Range("K77").Value = 18500
$ws.Range("L77").Value = 20250
$ws.Range("M77").Value = -13820
$ws.Range("N77").Value = -29610

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3048.4849
$ws.Range("I79").Value = 3027.5862
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3027.5862
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -1935.5862
$ws.Range("N79").Value = -5384

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2384.8948
$ws.Range("I98").Value = 1769.1875
$ws.Range("J98").Value = 5668.6665
$ws.Range("K98").Value = 1769.1875
$ws.Range("L98").Value = 5668.6665
$ws.Range("M98").Value = -271.1875
$ws.Range("N98").Value = -8664.666499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4569.8696
$ws.Range("I116").Value = 2550
$ws.Range("J116").Value = 6123.615
$ws.Range("K116").Value = 2550
$ws.Range("L116").Value = 6123.615
$ws.Range("M116").Value = 892
$ws.Range("N116").Value = -13007.615

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2384.8948
$ws.Range("I122").Value = 1769.1875
$ws.Range("J122").Value = 5668.6665
$ws.Range("K122").Value = 5307.5625
$ws.Range("L122").Value = 17005.9995
$ws.Range("M122").Value = -2857.5625
$ws.Range("N122").Value = -21905.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2339.1738
$ws.Range("I137").Value = 2084.5557
$ws.Range("K137").Value = 6253.6671
$ws.Range("M137").Value = -3703.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3363.318
$ws.Range("I63").Value = 2841.1667
$ws.Range("J63").Value = 3989.9
$ws.Range("K63").Value = 2841.1667
$ws.Range("L63").Value = 3989.9
$ws.Range("M63").Value = -2155.1667
$ws.Range("N63").Value = -5361.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3363.318
$ws.Range("I66").Value = 2841.1667
$ws.Range("J66").Value = 3989.9
$ws.Range("K66").Value = 14205.8335
$ws.Range("L66").Value = 19949.5
$ws.Range("M66").Value = -10773.8335
$ws.Range("N66").Value = -26813.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2228.513
$ws.Range("I132").Value = 1747.5333
$ws.Range("J132").Value = 3831.7778
$ws.Range("K132").Value = 5242.5999
$ws.Range("L132").Value = 11495.3334
$ws.Range("M132").Value = -2712.5999
$ws.Range("N132").Value = -16555.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2497.9167
$ws.Range("I105").Value = 2083.3333
$ws.Range("J105").Value = 2557.1428
$ws.Range("K105").Value = 2083.3333
$ws.Range("L105").Value = 2557.1428
$ws.Range("M105").Value = -336.3332999999998
$ws.Range("N105").Value = -6051.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2178.182
$ws.Range("I107").Value = 2384.8333
$ws.Range("J107").Value = 1930.2
$ws.Range("K107").Value = 2384.8333
$ws.Range("L107").Value = 1930.2
$ws.Range("M107").Value = -464.8332999999998
$ws.Range("N107").Value = -5770.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 883.5
$ws.Range("I16").Value = 820.2
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 820.2
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -533.2
$ws.Range("N16").Value = -1774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1596.4736
$ws.Range("I31").Value = 1377.0625
$ws.Range("J31").Value = 2766.6667
$ws.Range("K31").Value = 1377.0625
$ws.Range("L31").Value = 2766.6667
$ws.Range("M31").Value = -1082.0625
$ws.Range("N31").Value = -3356.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1596.4736
$ws.Range("I34").Value = 1377.0625
$ws.Range("J34").Value = 2766.6667
$ws.Range("K34").Value = 1377.0625
$ws.Range("L34").Value = 2766.6667
$ws.Range("M34").Value = -1175.0625
$ws.Range("N34").Value = -3170.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1704.3529
$ws.Range("I107").Value = 805.3333
$ws.Range("J107").Value = 2194.7273
$ws.Range("K107").Value = 805.3333
$ws.Range("L107").Value = 2194.7273
$ws.Range("M107").Value = 1114.6667
$ws.Range("N107").Value = -6034.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 883.5
$ws.Range("I113").Value = 820.2
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 820.2
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1349.8
$ws.Range("N113").Value = -5540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 780.6667
$ws.Range("I113").Value = 728.4286
$ws.Range("J113").Value = 853.8
$ws.Range("K113").Value = 2185.2858
$ws.Range("L113").Value = 2561.4
$ws.Range("M113").Value = -15.28579999999965
$ws.Range("N113").Value = -6901.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 865.62
$ws.Range("J131").Value = 894.1183
$ws.Range("L131").Value = 2682.3549
$ws.Range("N131").Value = -12762.3549

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5346.5
$ws.Range("I70").Value = 5136.727
$ws.Range("J70").Value = 5808
$ws.Range("K70").Value = 5136.727
$ws.Range("L70").Value = 5808
$ws.Range("M70").Value = -4866.727
$ws.Range("N70").Value = -6348

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5346.5
$ws.Range("I73").Value = 5136.727
$ws.Range("J73").Value = 5808
$ws.Range("K73").Value = 5136.727
$ws.Range("L73").Value = 5808
$ws.Range("M73").Value = -4200.727
$ws.Range("N73").Value = -7680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3083.8
$ws.Range("I132").Value = 2751.8667
$ws.Range("J132").Value = 4079.6
$ws.Range("K132").Value = 8255.6001
$ws.Range("L132").Value = 12238.8
$ws.Range("M132").Value = -5725.6001
$ws.Range("N132").Value = -17298.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 864.4815
$ws.Range("I46").Value = 868.1
$ws.Range("J46").Value = 862.35297
$ws.Range("K46").Value = 868.1
$ws.Range("L46").Value = 862.35297
$ws.Range("M46").Value = -680.1
$ws.Range("N46").Value = -1238.35297

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2515.6667
$ws.Range("I61").Value = 2120.1428
$ws.Range("J61").Value = 3900
$ws.Range("K61").Value = 2120.1428
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -1918.1428
$ws.Range("N61").Value = -4304

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1719.8
$ws.Range("I100").Value = 1399.75
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1399.75
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -858.75
$ws.Range("N100").Value = -4082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2515.6667
$ws.Range("I113").Value = 2120.1428
$ws.Range("J113").Value = 3900
$ws.Range("K113").Value = 2120.1428
$ws.Range("L113").Value = 3900
$ws.Range("M113").Value = 49.85719999999992
$ws.Range("N113").Value = -8240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 45423.5
$ws.Range("I41").Value = 78158.5
$ws.Range("J41").Value = 12688.5
$ws.Range("K41").Value = 78158.5
$ws.Range("L41").Value = 12688.5
$ws.Range("M41").Value = -77768.5
$ws.Range("N41").Value = -13468.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 31850
$ws.Range("J115").Value = 31850
$ws.Range("L115").Value = 31850
$ws.Range("N115").Value = -34984

Write-Host "Applied all leve profit updates"